$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 2210.6924
$ws.Range("I11").Value = 2210.6924
$ws.Range("K11").Value = 2210.6924
$ws.Range("M11").Value = -2070.6924
$ws.Range("H40").Value = 21432654
$ws.Range("J40").Value = 33336142
$ws.Range("L40").Value = 33336142
$ws.Range("N40").Value = -33336492
$ws.Range("H41").Value = 871.63635
$ws.Range("I41").Value = 467
$ws.Range("K41").Value = 467
$ws.Range("M41").Value = -27
$ws.Range("H43").Value = 6225.5713
$ws.Range("I43").Value = 2993.3333
$ws.Range("K43").Value = 2993.3333
$ws.Range("M43").Value = -2924.3333
$ws.Range("H47").Value = 8683
$ws.Range("I47").Value = 8683
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 8683
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -7711
$ws.Range("N47").ClearContents()
$ws.Range("H58").Value = 421.6
$ws.Range("I58").Value = 421.6
$ws.Range("K58").Value = 1264.8
$ws.Range("M58").Value = -1114.8
$ws.Range("H69").Value = 9885.888999999999
$ws.Range("I69").Value = 9885.888999999999
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 29657.667
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -28783.667
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 9885.888999999999
$ws.Range("I72").Value = 9885.888999999999
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 88973.00099999999
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -84605.00099999999
$ws.Range("N72").ClearContents()
$ws.Range("H86").Value = 3380.6
$ws.Range("J86").Value = 3425.25
$ws.Range("L86").Value = 3425.25
$ws.Range("N86").Value = -5671.25
$ws.Range("H89").Value = 3380.6
$ws.Range("J89").Value = 3425.25
$ws.Range("L89").Value = 17126.25
$ws.Range("N89").Value = -28358.25
$ws.Range("H98").Value = 6310426
$ws.Range("I98").Value = 8267711.5
$ws.Range("J98").Value = 2004399
$ws.Range("K98").Value = 8267711.5
$ws.Range("L98").Value = 2004399
$ws.Range("M98").Value = -8266213.5
$ws.Range("N98").Value = -2007395
$ws.Range("H122").Value = 6310426
$ws.Range("I122").Value = 8267711.5
$ws.Range("J122").Value = 2004399
$ws.Range("K122").Value = 24803134.5
$ws.Range("L122").Value = 6013197
$ws.Range("M122").Value = -24800684.5
$ws.Range("N122").Value = -6018097
$ws.Range("H132").Value = 3230
$ws.Range("I132").Value = 2924.2856
$ws.Range("J132").Value = 4300
$ws.Range("K132").Value = 8772.856800000001
$ws.Range("L132").Value = 12900
$ws.Range("M132").Value = -6242.856800000001
$ws.Range("N132").Value = -17960
$ws.Range("H141").Value = 2496.2354
$ws.Range("I141").Value = 2516.1724
$ws.Range("J141").Value = 2380.6
$ws.Range("K141").Value = 7548.5172
$ws.Range("L141").Value = 7141.799999999999
$ws.Range("M141").Value = -2368.5172
$ws.Range("N141").Value = -17501.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3024.1667
$ws.Range("I61").Value = 1662.4706
$ws.Range("K61").Value = 1662.4706
$ws.Range("M61").Value = -1450.4706
$ws.Range("H122").Value = 3365.125
$ws.Range("I122").Value = 2992.9473
$ws.Range("K122").Value = 8978.841899999999
$ws.Range("M122").Value = -6528.841899999999
$ws.Range("H136").Value = 3024.1667
$ws.Range("I136").Value = 1662.4706
$ws.Range("K136").Value = 4987.4118
$ws.Range("M136").Value = -2437.4118
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 681.3
$ws.Range("I94").Value = 607.9375
$ws.Range("J94").Value = 974.75
$ws.Range("K94").Value = 607.9375
$ws.Range("L94").Value = 974.75
$ws.Range("M94").Value = -156.9375
$ws.Range("N94").Value = -1876.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2316.0557
$ws.Range("I16").Value = 2285.7
$ws.Range("J16").Value = 2354
$ws.Range("K16").Value = 2285.7
$ws.Range("L16").Value = 2354
$ws.Range("M16").Value = -1998.7
$ws.Range("N16").Value = -2928
$ws.Range("H22").Value = 331.33334
$ws.Range("I22").Value = 338.43478
$ws.Range("J22").Value = 168
$ws.Range("K22").Value = 338.43478
$ws.Range("L22").Value = 168
$ws.Range("M22").Value = 11.56522000000001
$ws.Range("N22").Value = -868
$ws.Range("H113").Value = 2316.0557
$ws.Range("I113").Value = 2285.7
$ws.Range("J113").Value = 2354
$ws.Range("K113").Value = 2285.7
$ws.Range("L113").Value = 2354
$ws.Range("M113").Value = -115.6999999999998
$ws.Range("N113").Value = -6694
$ws.Range("H122").Value = 5947.6
$ws.Range("I122").Value = 5742.143
$ws.Range("K122").Value = 17226.429
$ws.Range("M122").Value = -14776.429
$ws.Range("H132").Value = 3023.0278
$ws.Range("I132").Value = 2898.8696
$ws.Range("J132").Value = 3242.6924
$ws.Range("K132").Value = 8696.6088
$ws.Range("L132").Value = 9728.0772
$ws.Range("M132").Value = -6166.6088
$ws.Range("N132").Value = -14788.0772
$ws.Range("H134").Value = 5271.2
$ws.Range("J134").Value = 3381.4285
$ws.Range("L134").Value = 10144.2855
$ws.Range("N134").Value = -15214.2855
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 823.5
$ws.Range("I9").Value = 306
$ws.Range("J9").Value = 1045.2858
$ws.Range("K9").Value = 918
$ws.Range("L9").Value = 3135.8574
$ws.Range("M9").Value = -694
$ws.Range("N9").Value = -3583.8574
$ws.Range("H37").Value = 58879732
$ws.Range("J37").Value = 58879732
$ws.Range("L37").Value = 176639196
$ws.Range("N37").Value = -176639420
$ws.Range("H68").Value = 11905859
$ws.Range("I68").Value = 27778306
$ws.Range("J68").Value = 1524.75
$ws.Range("K68").Value = 83334918
$ws.Range("L68").Value = 4574.25
$ws.Range("M68").Value = -83334107
$ws.Range("N68").Value = -6196.25
$ws.Range("H71").Value = 11905859
$ws.Range("I71").Value = 27778306
$ws.Range("J71").Value = 1524.75
$ws.Range("K71").Value = 250004754
$ws.Range("L71").Value = 13722.75
$ws.Range("M71").Value = -250000698
$ws.Range("N71").Value = -21834.75
$ws.Range("H98").Value = 479.22223
$ws.Range("I98").Value = 316.25
$ws.Range("J98").Value = 609.6
$ws.Range("K98").Value = 948.75
$ws.Range("L98").Value = 1828.8
$ws.Range("M98").Value = 549.25
$ws.Range("N98").Value = -4824.8
$ws.Range("H131").Value = 17242864
$ws.Range("I131").Value = 125000910
$ws.Range("K131").Value = 375002730
$ws.Range("M131").Value = -374997690
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3850.739
$ws.Range("I97").Value = 3319.7058
$ws.Range("J97").Value = 5355.3335
$ws.Range("K97").Value = 3319.7058
$ws.Range("L97").Value = 5355.3335
$ws.Range("M97").Value = -2823.7058
$ws.Range("N97").Value = -6347.3335
$ws.Range("H107").Value = 2799.875
$ws.Range("J107").Value = 2725.5
$ws.Range("L107").Value = 2725.5
$ws.Range("N107").Value = -6565.5
$ws.Range("H132").Value = 31790.107
$ws.Range("I132").Value = 35129.332
$ws.Range("J132").Value = 4241.5
$ws.Range("K132").Value = 105387.996
$ws.Range("L132").Value = 12724.5
$ws.Range("M132").Value = -102857.996
$ws.Range("N132").Value = -17784.5
$ws.Range("H136").Value = 33666.668
$ws.Range("J136").Value = 33666.668
$ws.Range("L136").Value = 101000.004
$ws.Range("N136").Value = -106100.004
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 19067.75
$ws.Range("J43").Value = 19067
$ws.Range("L43").Value = 19067
$ws.Range("N43").Value = -19453
$ws.Range("H53").Value = 52749.5
$ws.Range("I53").Value = 5500
$ws.Range("J53").Value = 99999
$ws.Range("K53").Value = 5500
$ws.Range("L53").Value = 99999
$ws.Range("M53").Value = -4982
$ws.Range("N53").Value = -101035
$ws.Range("H68").Value = 9976.117
$ws.Range("I68").Value = 11072.8
$ws.Range("J68").Value = 1751
$ws.Range("K68").Value = 11072.8
$ws.Range("L68").Value = 1751
$ws.Range("M68").Value = -10323.8
$ws.Range("N68").Value = -3249
$ws.Range("H71").Value = 9976.117
$ws.Range("I71").Value = 11072.8
$ws.Range("J71").Value = 1751
$ws.Range("K71").Value = 55364
$ws.Range("L71").Value = 8755
$ws.Range("M71").Value = -51620
$ws.Range("N71").Value = -16243
$ws.Range("H122").Value = 3690.077
$ws.Range("I122").Value = 4174.4
$ws.Range("K122").Value = 12523.2
$ws.Range("M122").Value = -10073.2
$ws.Range("H132").Value = 4893.6
$ws.Range("I132").Value = 5219.5713
$ws.Range("K132").Value = 15658.7139
$ws.Range("M132").Value = -13128.7139
$ws.Range("H136").Value = 1294.85
$ws.Range("I136").Value = 1194.8667
$ws.Range("J136").Value = 1594.8
$ws.Range("K136").Value = 3584.6001
$ws.Range("L136").Value = 4784.4
$ws.Range("M136").Value = -1034.6001
$ws.Range("N136").Value = -9884.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9097.200000000001
$ws.Range("J62").Value = 9496.5
$ws.Range("L62").Value = 9496.5
$ws.Range("N62").Value = -10744.5
$ws.Range("H65").Value = 9097.200000000001
$ws.Range("J65").Value = 9496.5
$ws.Range("L65").Value = 47482.5
$ws.Range("N65").Value = -53722.5
$ws.Range("H112").Value = 67127.836
$ws.Range("J112").Value = 67127.836
$ws.Range("L112").Value = 67127.836
$ws.Range("N112").Value = -70081.836
$ws.Range("H122").Value = 10259.167
$ws.Range("I122").Value = 2529.8147
$ws.Range("K122").Value = 7589.4441
$ws.Range("M122").Value = -5139.4441
$ws.Range("H132").Value = 1756.1428
$ws.Range("I132").Value = 1698.9231
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5096.7693
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -2566.7693
$ws.Range("N132").Value = -12560
